$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Insert a new row above current row 2 (shifts existing rows 2..9 down to 3..10)
$ws.Rows.Item(2).Insert()

# Populate the new filter row
$ws.Range("A2").Value = "TOC_Filter"
$ws.Range("B2:G2").Value = "All TOCs"

# Select B10, matching the end state in the saved file
$ws.Range("B10").Select()
